$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=-0.1541247462703282; C=0.852839648487878;  D=1.456859763062113;  E=1.207004458592475; F=1.209035744782678; G=51},
    @{Row=3;  B=0.1946594554508169;  C=0.9013317542106349; D=1.740862592224634;  E=1.319417520053692; F=1.318227895439372; G=50},
    @{Row=4;  B=0.09587990046027713;C=0.7813449169389153; D=1.273947737936718;  E=1.128692933413122; F=1.136267503977733; G=49},
    @{Row=5;  B=0.2010189561291209; C=0.8357504963916088; D=1.399781063783387;  E=1.183123435565109; F=1.1782594200759;   G=48},
    @{Row=6;  B=0.1134028742995598; C=0.7037266954146911; D=0.9458163126725945; E=0.9725308800611909;F=0.976338955044726; G=47},
    @{Row=7;  B=0.1165933518286177; C=0.7311150576301884; D=1.063932431078233;  E=1.03147100350821;  F=1.038617324277796; G=38},
    @{Row=8;  B=0.1240158675361395; C=0.6673632257439414; D=0.9510627036727184; E=0.9752244375899931;F=0.9806497622900503;G=37},
    @{Row=9;  B=-0.06833470819847304;C=0.3931001917406148;D=0.2556107661499554; E=0.5055796338362092;F=0.5139538519170813;G=20},
    @{Row=10; B=-0.001481533045873319;C=0.394696832903608;D=0.2639446156350248; E=0.5137554044825463;F=0.5347313553733598;G=13},
    @{Row=11; B=0.1038035452512703; C=0.3551067294956168; D=0.1935206269752047; E=0.4399097941342119;F=0.4779454087137884;G=5}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
